## CSV-PWM-Table.xlsx — Version 13 edits
## - PWM freq decreased to 18KHz: Counter Max (H2) 210 -> 444, derived cells updated
## - New "PWM Freq." column (replaces "Check") with Hz number format
## - New helper cells L2/M2
## - Header row + counter-max cell get a bottom-thick rule / yellow highlight / box border
## - Column I widened, dimension/selection follow the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): "Check" becomes "PWM Freq.", Scale/Counter Max
#    keep their styles (reindexed), PWM Freq. reuses the old Counter-Max-style
#    border/center/bold look.
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = "PWM Freq."

# thin box border + bold centered font, same look previously used for H1/I1
foreach ($addr in @("G1","H1","I1")) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108   # xlCenter
}
$ws.Range("H1:I1").Font.Bold = $true
$ws.Range("G1:I1").Borders.LineStyle = 1
$ws.Range("G1:I1").Borders.Weight = 2

# Row 1 gets a heavier bottom rule (thick bottom border) to set the header off
$ws.Range("A1:M1").Borders.Item(9).LineStyle = 1
$ws.Range("A1:M1").Borders.Item(9).Weight = 3
$ws.Rows.Item(1).RowHeight = 15

# ---------------------------------------------------------------------------
# 2. Row 2 — the PWM parameter block.
# ---------------------------------------------------------------------------
# Counter Max input changes from 210 -> 444
$ws.Range("H2").Value = 444

# Scale formula changes (was ROUND(H2*256/254,0))
$ws.Range("G2").Formula = "=H2/2-7"

# "Check" column repurposed to show the resulting PWM frequency
$ws.Range("I2").Formula = "=16000000/(H2*2)"
$ws.Range("I2").NumberFormat = "0 ""Hz"""

# New helper cells (not previously present)
$ws.Range("L2").Formula = "=TRUNC(254*G2/256)"
$ws.Range("M2").Formula = "=H2/2-L2"

# Counter Max cell: yellow fill + medium (thick) box border, red bold font
$ws.Range("H2").Interior.Color = 65535
$ws.Range("H2").Font.Color = 255
$ws.Range("H2").Borders.LineStyle = 1
$ws.Range("H2").Borders.Weight = -4138

# PWM Freq. result cell: thin border on right/top/bottom (matches the header box)
$ws.Range("I2").Borders.Item(10).LineStyle = 1   # right
$ws.Range("I2").Borders.Item(8).LineStyle = 1    # top
$ws.Range("I2").Borders.Item(9).LineStyle = 1    # bottom

$ws.Rows.Item(2).RowHeight = 15

# ---------------------------------------------------------------------------
# 3. Column I is no longer auto-fit — give it an explicit width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 11.77734375

# ---------------------------------------------------------------------------
# 4. Selection follows the parameter cell that now matters most.
# ---------------------------------------------------------------------------
$ws.Range("G2").Select()

$wb.Application.Calculate()
